$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 44
$ws.Range("H44").Value = 50000  # was 45600
$ws.Range("J44").Value = 50000  # was 45600
$ws.Range("L44").Value = 50000  # was 45600
$ws.Range("N44").Value = -50924  # was -46524
# Row 126
$ws.Range("H126").Value = 41926  # was 41901.668
$ws.Range("J126").Value = 41926  # was 41901.668
$ws.Range("L126").Value = 41926  # was 41901.668
$ws.Range("N126").Value = -51806  # was -51781.668
# Row 132
$ws.Range("H132").Value = 21045010  # was 21959962
$ws.Range("I132").Value = 21742294  # was 22730528
$ws.Range("K132").Value = 65226882  # was 68191584
$ws.Range("M132").Value = -65224352  # was -68189054
# Row 133
$ws.Range("H133").Value = 41043.332  # was 38251.3
$ws.Range("J133").Value = 41043.332  # was 38251.3
$ws.Range("L133").Value = 41043.332  # was 38251.3
$ws.Range("N133").Value = -51163.332  # was -48371.3
# Row 134
$ws.Range("H134").Value = 52982.38  # was 51860.527
$ws.Range("J134").Value = 52982.38  # was 51860.527
$ws.Range("L134").Value = 52982.38  # was 51860.527
$ws.Range("N134").Value = -63122.38  # was -62000.527
# Row 136
$ws.Range("H136").Value = 54333.332  # was 45292.5
$ws.Range("J136").Value = 54333.332  # was 45292.5
$ws.Range("L136").Value = 54333.332  # was 45292.5
$ws.Range("N136").Value = -64533.332  # was -55492.5
# Row 139
$ws.Range("H139").Value = 45780  # was 38360
$ws.Range("J139").Value = 45780  # was 38360
$ws.Range("L139").Value = 45780  # was 38360
$ws.Range("N139").Value = -56060  # was -48640
# Row 140
$ws.Range("H140").Value = 62254.547  # was 46318.387
$ws.Range("J140").Value = 62254.547  # was 46318.387
$ws.Range("L140").Value = 62254.547  # was 46318.387
$ws.Range("N140").Value = -72614.54699999999  # was -56678.387
# Row 141
$ws.Range("H141").Value = 9844.691999999999  # was 10415.083
$ws.Range("I141").Value = 11818.1  # was 12797.889
$ws.Range("K141").Value = 35454.3  # was 38393.667
$ws.Range("M141").Value = -30274.3  # was -33213.667

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 573.25  # was 559.5
$ws.Range("I2").Value = 455.5  # was 424.4
$ws.Range("J2").Value = 691  # was 1235
$ws.Range("K2").Value = 455.5  # was 424.4
$ws.Range("L2").Value = 691  # was 1235
$ws.Range("M2").Value = -342.5  # was -311.4
$ws.Range("N2").Value = -917  # was -1461
# Row 110
$ws.Range("H110").Value = 902.75  # was 522.5333000000001
$ws.Range("I110").Value = 720.5  # was 382.30768
$ws.Range("J110").Value = 1085  # was 1434
$ws.Range("K110").Value = 720.5  # was 382.30768
$ws.Range("L110").Value = 1085  # was 1434
$ws.Range("M110").Value = 1324.5  # was 1662.69232
$ws.Range("N110").Value = -5175  # was -5524
# Row 116
$ws.Range("H116").Value = 573.25  # was 559.5
$ws.Range("I116").Value = 455.5  # was 424.4
$ws.Range("J116").Value = 691  # was 1235
$ws.Range("K116").Value = 455.5  # was 424.4
$ws.Range("L116").Value = 691  # was 1235
$ws.Range("M116").Value = 1838.5  # was 1869.6
$ws.Range("N116").Value = -5279  # was -5823

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 573.25  # was 559.5
$ws.Range("I3").Value = 455.5  # was 424.4
$ws.Range("J3").Value = 691  # was 1235
$ws.Range("K3").Value = 455.5  # was 424.4
$ws.Range("L3").Value = 691  # was 1235
$ws.Range("M3").Value = -341.5  # was -310.4
$ws.Range("N3").Value = -919  # was -1463

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7054.2617  # was 7156.775
$ws.Range("I31").Value = 3163.9  # was 3163.5
$ws.Range("J31").Value = 10590.954  # was 11150.05
$ws.Range("K31").Value = 3163.9  # was 3163.5
$ws.Range("L31").Value = 10590.954  # was 11150.05
$ws.Range("M31").Value = -2868.9  # was -2868.5
$ws.Range("N31").Value = -11180.954  # was -11740.05
# Row 34
$ws.Range("H34").Value = 7054.2617  # was 7156.775
$ws.Range("I34").Value = 3163.9  # was 3163.5
$ws.Range("J34").Value = 10590.954  # was 11150.05
$ws.Range("K34").Value = 3163.9  # was 3163.5
$ws.Range("L34").Value = 10590.954  # was 11150.05
$ws.Range("M34").Value = -2961.9  # was -2961.5
$ws.Range("N34").Value = -10994.954  # was -11554.05
# Row 36
$ws.Range("H36").Value = 1950  # was 1500
$ws.Range("I36").Value = 1950  # was 1500
$ws.Range("K36").Value = 1950  # was 1500
$ws.Range("M36").Value = -1562  # was -1112
# Row 40
$ws.Range("H40").Value = 1950  # was 1500
$ws.Range("I40").Value = 1950  # was 1500
$ws.Range("K40").Value = 1950  # was 1500
$ws.Range("M40").Value = -1790  # was -1340
# Row 86
$ws.Range("H86").Value = 3233.2222  # was 4116.5
$ws.Range("I86").Value = 2550  # was 4450
$ws.Range("J86").Value = 4599.6665  # was 3949.75
$ws.Range("K86").Value = 2550  # was 4450
$ws.Range("L86").Value = 4599.6665  # was 3949.75
$ws.Range("M86").Value = -1427  # was -3327
$ws.Range("N86").Value = -6845.6665  # was -6195.75
# Row 89
$ws.Range("H89").Value = 3233.2222  # was 4116.5
$ws.Range("I89").Value = 2550  # was 4450
$ws.Range("J89").Value = 4599.6665  # was 3949.75
$ws.Range("K89").Value = 12750  # was 22250
$ws.Range("L89").Value = 22998.3325  # was 19748.75
$ws.Range("M89").Value = -7134  # was -16634
$ws.Range("N89").Value = -34230.3325  # was -30980.75
# Row 98
$ws.Range("H98").Value = 42676.5  # was 50000
$ws.Range("J98").Value = 42676.5  # was 50000
$ws.Range("L98").Value = 42676.5  # was 50000
$ws.Range("N98").Value = -47168.5  # was -54492
# Row 124
$ws.Range("H124").Value = 47999  # was 48000
$ws.Range("J124").Value = 47999  # was 48000
$ws.Range("L124").Value = 47999  # was 48000
$ws.Range("N124").Value = -52909  # was -52910
# Row 127
$ws.Range("H127").Value = 43030  # was 42873.75
$ws.Range("J127").Value = 43030  # was 42873.75
$ws.Range("L127").Value = 43030  # was 42873.75
$ws.Range("N127").Value = -52950  # was -52793.75
# Row 132
$ws.Range("H132").Value = 5690.375  # was 2695.7
$ws.Range("I132").Value = 3705  # was 1373.4667
$ws.Range("J132").Value = 8999.333000000001  # was 6662.4
$ws.Range("K132").Value = 11115  # was 4120.4001
$ws.Range("L132").Value = 26997.999  # was 19987.2
$ws.Range("M132").Value = -8585  # was -1590.4001
$ws.Range("N132").Value = -32057.999  # was -25047.2
# Row 134
$ws.Range("H134").Value = 8453.588  # was 5721.2593
$ws.Range("I134").Value = 10700.272  # was 6945.8887
$ws.Range("J134").Value = 4334.6665  # was 3272
$ws.Range("K134").Value = 32100.816  # was 20837.6661
$ws.Range("L134").Value = 13003.9995  # was 9816
$ws.Range("M134").Value = -29565.816  # was -18302.6661
$ws.Range("N134").Value = -18073.9995  # was -14886

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 9999.333000000001  # was 9599.6
$ws.Range("I74").Value = 10000  # was 9333.333000000001
$ws.Range("K74").Value = 30000  # was 27999.999
$ws.Range("M74").Value = -28939  # was -26938.999
# Row 77
$ws.Range("H77").Value = 9999.333000000001  # was 9599.6
$ws.Range("I77").Value = 10000  # was 9333.333000000001
$ws.Range("K77").Value = 90000  # was 83999.997
$ws.Range("M77").Value = -84696  # was -78695.997
# Row 132
$ws.Range("H132").Value = 2171.6296  # was 2134.2666
$ws.Range("I132").Value = 961.46155  # was 999.375
$ws.Range("J132").Value = 3295.3572  # was 2546.9546
$ws.Range("K132").Value = 8653.15395  # was 8994.375
$ws.Range("L132").Value = 29658.2148  # was 22922.5914
$ws.Range("M132").Value = -6123.15395  # was -6464.375
$ws.Range("N132").Value = -34718.2148  # was -27982.5914
# Row 133
$ws.Range("H133").Value = 4646.857  # was 4960.6665
$ws.Range("I133").Value = 4646.857  # was 4992.8
$ws.Range("J133").Value = 0  # was 4800
$ws.Range("K133").Value = 13940.571  # was 14978.4
$ws.Range("L133").Value = 0  # was 14400
$ws.Range("M133").Value = -8880.571  # was -9918.400000000001
$ws.Range("N133").ClearContents()  # was -24520
# Row 134
$ws.Range("H134").Value = 4207.1  # was 4412.6665
$ws.Range("I134").Value = 3619.5625  # was 3877.1428
$ws.Range("J134").Value = 4878.5713  # was 4881.25
$ws.Range("K134").Value = 10858.6875  # was 11631.4284
$ws.Range("L134").Value = 14635.7139  # was 14643.75
$ws.Range("M134").Value = -5788.6875  # was -6561.428400000001
$ws.Range("N134").Value = -24775.7139  # was -24783.75
# Row 138
$ws.Range("H138").Value = 2596  # was 2211.6667
$ws.Range("I138").Value = 490  # was 923.3333
$ws.Range("J138").Value = 4000  # was 3500
$ws.Range("K138").Value = 1470  # was 2769.9999
$ws.Range("L138").Value = 12000  # was 10500
$ws.Range("M138").Value = 3670  # was 2370.0001
$ws.Range("N138").Value = -22280  # was -20780
# Row 139
$ws.Range("H139").Value = 1056.2084  # was 1535.45
$ws.Range("I139").Value = 884.73914  # was 1044.3125
$ws.Range("J139").Value = 5000  # was 3500
$ws.Range("K139").Value = 2654.21742  # was 3132.9375
$ws.Range("L139").Value = 15000  # was 10500
$ws.Range("M139").Value = 2485.78258  # was 2007.0625
$ws.Range("N139").Value = -25280  # was -20780
# Row 140
$ws.Range("H140").Value = 3957.6924  # was 6021.6665
$ws.Range("I140").Value = 8237.5  # was 10743.333
$ws.Range("J140").Value = 2055.5557  # was 1300
$ws.Range("K140").Value = 24712.5  # was 32229.999
$ws.Range("L140").Value = 6166.6671  # was 3900
$ws.Range("M140").Value = -19532.5  # was -27049.999
$ws.Range("N140").Value = -16526.6671  # was -14260
# Row 141
$ws.Range("H141").Value = 7954.6665  # was 8656.429
$ws.Range("I141").Value = 8108.5713  # was 9132.5
$ws.Range("K141").Value = 24325.7139  # was 27397.5
$ws.Range("M141").Value = -19145.7139  # was -22217.5

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 103
$ws.Range("H103").Value = 21786.666  # was 25944
$ws.Range("J103").Value = 21786.666  # was 25944
$ws.Range("L103").Value = 21786.666  # was 25944
$ws.Range("N103").Value = -24130.666  # was -28288
# Row 111
$ws.Range("H111").Value = 29420  # was 29520
$ws.Range("J111").Value = 29420  # was 29520
$ws.Range("L111").Value = 29420  # was 29520
$ws.Range("N111").Value = -35554  # was -35654
# Row 118
$ws.Range("H118").Value = 25195.715  # was 25532.5
$ws.Range("J118").Value = 25195.715  # was 25532.5
$ws.Range("L118").Value = 25195.715  # was 25532.5
$ws.Range("N118").Value = -28509.715  # was -28846.5
# Row 122
$ws.Range("H122").Value = 4200.2  # was 4218.4116
$ws.Range("I122").Value = 3610.2727  # was 3476.0833
$ws.Range("J122").Value = 5822.5  # was 6000
$ws.Range("K122").Value = 10830.8181  # was 10428.2499
$ws.Range("L122").Value = 17467.5  # was 18000
$ws.Range("M122").Value = -8380.8181  # was -7978.249899999999
$ws.Range("N122").Value = -22367.5  # was -22900

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 115
$ws.Range("H115").Value = 39000  # was 39800
$ws.Range("J115").Value = 39000  # was 39800
$ws.Range("L115").Value = 39000  # was 39800
$ws.Range("N115").Value = -41350  # was -42150
# Row 119
$ws.Range("H119").Value = 34860  # was 34920
$ws.Range("J119").Value = 34860  # was 34920
$ws.Range("L119").Value = 34860  # was 34920
$ws.Range("N119").Value = -44536  # was -44596
# Row 121
$ws.Range("H121").Value = 27661.428  # was 26631.25
$ws.Range("J121").Value = 27661.428  # was 26631.25
$ws.Range("L121").Value = 27661.428  # was 26631.25
$ws.Range("N121").Value = -31155.428  # was -30125.25
# Row 132
$ws.Range("H132").Value = 4932.952  # was 5050.1577
$ws.Range("I132").Value = 3799.6875  # was 3910.3333
$ws.Range("J132").Value = 8559.4  # was 9324.5
$ws.Range("K132").Value = 11399.0625  # was 11730.9999
$ws.Range("L132").Value = 25678.2  # was 27973.5
$ws.Range("M132").Value = -8869.0625  # was -9200.999899999999
$ws.Range("N132").Value = -30738.2  # was -33033.5

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 19610538  # was 19610426
$ws.Range("I132").Value = 1483.1666  # was 1324.875
$ws.Range("J132").Value = 30306386  # was 37040736
$ws.Range("K132").Value = 4449.4998  # was 3974.625
$ws.Range("L132").Value = 90919158  # was 111122208
$ws.Range("M132").Value = -1919.4998  # was -1444.625
$ws.Range("N132").Value = -90924218  # was -111127268
